{"js": "// Office.js (Word JavaScript API) edit script.\n//\n// This reproduces the change described by the diff:\n//  1. Remove the \"Meta description: ...\" paragraph that follows the\n//     \"Play Book of Time Slot Game Free - Review\" Heading1 paragraph\n//     at the top of the document.\n//  2. At the very end of the document, just before the final \"Prompt: ...\"\n//     paragraph:\n//       a. insert a new bold paragraph with the text\n//          \"Play Book of Time Slot Game Free - Review\"\n//       b. replace the text of the final paragraph (still italic) with\n//          \"Read our review of Book of Time online slot game. Play Book\n//          of Time for free and learn about its features, bonuses, and\n//          design.\" (i.e. drop the leading \"Prompt: Create a feature\n//          image...\" content and the \"Meta description: \" label).\n\nconst body = context.document.body;\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\nconst paragraphs = body.paragraphs.items;\n\n// --- Step 1: delete the \"Meta description\" paragraph near the top ------\nconst metaPara = paragraphs.find((p) => p.text.indexOf(\"Meta description\") === 0);\nif (metaPara) {\n  metaPara.delete();\n}\n\n// --- Step 2: rework the final paragraph (was the \"Prompt: ...\" text) ---\n// Re-load paragraphs since the collection changed after the delete above.\nbody.paragraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = body.paragraphs.items;\nconst lastParagraph = items[items.length - 1];\n\n// 2a. Insert a new bold heading-like paragraph right before it.\nconst newHeadingParagraph = lastParagraph.insertParagraph(\n  \"Play Book of Time Slot Game Free - Review\",\n  \"Before\"\n);\nnewHeadingParagraph.font.bold = true;\nnewHeadingParagraph.font.italic = false;\n\n// 2b. Swap out the prompt text for the meta-description text, keeping\n// the paragraph's existing italic formatting.\nlastParagraph.insertText(\n  \"Read our review of Book of Time online slot game. Play Book of Time for free and learn about its features, bonuses, and design.\",\n  \"Replace\"\n);\n\nawait context.sync();\n", "ps1": "# Word COM interop script.\n#\n# Reproduces the change described by the diff:\n#  1. Remove the \"Meta description: ...\" paragraph that follows the\n#     \"Play Book of Time Slot Game Free - Review\" Heading1 paragraph\n#     at the top of the document.\n#  2. At the very end of the document, just before the final \"Prompt: ...\"\n#     paragraph:\n#       a. insert a new bold paragraph with the text\n#          \"Play Book of Time Slot Game Free - Review\"\n#       b. replace the text of the final paragraph (still italic) with\n#          \"Read our review of Book of Time online slot game. Play Book\n#          of Time for free and learn about its features, bonuses, and\n#          design.\" (i.e. drop the leading \"Prompt: Create a feature\n#          image...\" content and the \"Meta description: \" label).\n\n$d = $word.ActiveDocument\n\n# --- Step 1: delete the \"Meta description\" paragraph near the top ------\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $para = $d.Paragraphs.Item($i)\n    if ($para.Range.Text.IndexOf(\"Meta description\") -eq 0) {\n        $para.Range.Delete()\n        break\n    }\n}\n\n# --- Step 2: rework the final paragraph (was the \"Prompt: ...\" text) ---\n$n = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($n)\n\n# 2a. Insert a new paragraph right before the final one and fill it with\n#     the bold heading text (explicitly not italic, since the\n#     neighbouring paragraph is italic and would otherwise be inherited).\n$lastPara.Range.InsertParagraphBefore()\n$newPara = $d.Paragraphs.Item($n)\n$newRange = $newPara.Range\n$newRange.MoveEnd(1, -1)\n$newRange.Text = \"Play Book of Time Slot Game Free - Review\"\n$newRange.Font.Bold = $true\n$newRange.Font.Italic = $false\n\n# 2b. Swap out the prompt text for the meta-description text, keeping\n#     the paragraph's existing (italic) formatting intact.\n$promptPara = $d.Paragraphs.Item($n + 1)\n$promptRange = $promptPara.Range\n$promptRange.MoveEnd(1, -1)\n$promptRange.Text = \"Read our review of Book of Time online slot game. Play Book of Time for free and learn about its features, bonuses, and design.\"\n"}
